$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Sending cluster = ECs, Target cluster = ECs)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dsc2"
$ws.Range("C2").Value = "Dsg1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3177946666666667
$ws.Range("H2").Value = 0.953384
$ws.Range("I2").Value = 0.2153608586026293
$ws.Range("J2").Value = 0.2153608586026293
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.004083333333333334
$ws.Range("N2").Value = 0.01225
$ws.Range("O2").Value = 0.1210581968752162
$ws.Range("P2").Value = 0.1210581968752162
$ws.Range("Q2").Value = 0.001297661555555556
$ws.Range("R2").Value = 0.011678954
$ws.Range("S2").Value = 0.02607119721993269
$ws.Range("T2").Value = 0.02607119721993269

# New row 3 (Sending cluster = ECs, Target cluster = MuSCs)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dsc2"
$ws.Range("C3").Value = "Dsg1a"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3177946666666667
$ws.Range("H3").Value = 0.953384
$ws.Range("I3").Value = 0.2153608586026293
$ws.Range("J3").Value = 0.2153608586026293
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.029647
$ws.Range("N3").Value = 0.08894100000000001
$ws.Range("O3").Value = 0.8789418031247839
$ws.Range("P3").Value = 0.8789418031247839
$ws.Range("Q3").Value = 0.009421658482666668
$ws.Range("R3").Value = 0.084794926344
$ws.Range("S3").Value = 0.1892896613826966
$ws.Range("T3").Value = 0.1892896613826966

# New row 4 (Sending cluster = MuSCs, Target cluster = ECs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Dsc2"
$ws.Range("C4").Value = "Dsg1a"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.157843333333333
$ws.Range("H4").Value = 3.47353
$ws.Range("I4").Value = 0.7846391413973707
$ws.Range("J4").Value = 0.7846391413973708
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.004083333333333334
$ws.Range("N4").Value = 0.01225
$ws.Range("O4").Value = 0.1210581968752162
$ws.Range("P4").Value = 0.1210581968752162
$ws.Range("Q4").Value = 0.004727860277777779
$ws.Range("R4").Value = 0.0425507425
$ws.Range("S4").Value = 0.09498699965528348
$ws.Range("T4").Value = 0.0949869996552835

# New row 5 (Sending cluster = MuSCs, Target cluster = MuSCs)
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Dsc2"
$ws.Range("C5").Value = "Dsg1a"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.157843333333333
$ws.Range("H5").Value = 3.47353
$ws.Range("I5").Value = 0.7846391413973707
$ws.Range("J5").Value = 0.7846391413973708
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.029647
$ws.Range("N5").Value = 0.08894100000000001
$ws.Range("O5").Value = 0.8789418031247839
$ws.Range("P5").Value = 0.8789418031247839
$ws.Range("Q5").Value = 0.03432658130333334
$ws.Range("R5").Value = 0.30893923173
$ws.Range("S5").Value = 0.6896521417420872
$ws.Range("T5").Value = 0.6896521417420873
